$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original Text formatting so that values such as
# "1.000" or "0.000007902" are not coerced into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.220.95'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.853.77'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.6982'
$ws.Range("E5").Value = '  +2.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.66'
$ws.Range("E6").Value = '  +0.12%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08113'
$ws.Range("E8").Value = '  +5.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3015'
$ws.Range("E9").Value = '  -0.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.46'
$ws.Range("E10").Value = '  +1.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08173'
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.869.10'
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.192'
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7048'
$ws.Range("E14").Value = '  -2.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.72'
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.237.43'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.815'
$ws.Range("E17").Value = '  +1.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007902'
$ws.Range("E18").Value = '  +1.36%  '
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '236.76'
$ws.Range("E20").Value = '  +1.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.104.86'
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.439'
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.90'
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.872'
$ws.Range("E26").Value = '  -0.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1413'
$ws.Range("E27").Value = '  -1.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.04'
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.914'
$ws.Range("E30").Value = '  +0.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.477'
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.353'
$ws.Range("E32").Value = '  -3.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.022'
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05183'
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.160'
$ws.Range("E35").Value = '  -1.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7191'
$ws.Range("E36").Value = '  +2.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9986'
$ws.Range("E37").Value = '  -2.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.688'
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01847'
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.718'
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9332'
$ws.Range("E41").Value = '  +2.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.145.26'
$ws.Range("E42").Value = '  +3.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.002'
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4248'
$ws.Range("E44").Value = '  -0.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.19'
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.84'
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("E48").Value = '  -3.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.750'
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.004.11'
$ws.Range("E50").Value = '  -0.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.145'
$ws.Range("E51").Value = '  +0.15%  '
